# Self-Evaluation-Protocol.xlsx — score corrections + a couple of ticked-off
# checklist cells, plus the window's scroll/selection state left where the
# author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "AngularJS Project Structure" score bumped from 9 -> 10
$ws.Range("C19").Value = 10

# "Edit User Profile" score was blank, now filled in with 10
$ws.Range("C26").Value = 10

# C44 is =SUM(C6:C43) and recalculates automatically (299 -> 310)

# Window state: scrolled further down and a different cell selected
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C20").Select()
